# Update the "想去人数" (number of people interested) figures for several
# events that appear on both the "展览" (exhibitions) sheet and the
# "全部类型" (all types) summary sheet.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsAll  = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet ---
$wsExpo.Range("F2").Value = 96
$wsExpo.Range("F3").Value = 4060
$wsExpo.Range("F4").Value = 2371
$wsExpo.Range("F8").Value = 29
$wsExpo.Range("F11").Value = 82
$wsExpo.Range("F12").Value = 135
$wsExpo.Range("F13").Value = 1521
$wsExpo.Range("F14").Value = 271
$wsExpo.Range("F15").Value = 2899
$wsExpo.Range("F16").Value = 200

# --- 全部类型 sheet ---
$wsAll.Range("F2").Value = 96
$wsAll.Range("F3").Value = 4060
$wsAll.Range("F4").Value = 2371
$wsAll.Range("F9").Value = 29
$wsAll.Range("F13").Value = 82
$wsAll.Range("F14").Value = 135
$wsAll.Range("F17").Value = 1521
$wsAll.Range("F18").Value = 271
$wsAll.Range("F19").Value = 2899
$wsAll.Range("F20").Value = 200
